# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# 1) Update the "last updated" timestamp string (cell A1)
$ws.Range("A1").Value = "Datos actualizados a 5 de Agosto de 2020 a las 09:22"

# 2) Ucrania overtakes Republica Dominicana in the ranking (row 37/38 swap
#    positions; Ucrania gets fresh numbers, Republica Dominicana keeps its
#    previous numbers but drops one spot).
$ws.Range("A37").Value = "Ucrania"
$ws.Range("B37").Value = 75490
$ws.Range("C37").Value = 1271
$ws.Range("D37").Value = 41527
$ws.Range("E37").Value = 32175
$ws.Range("F37").Value = 0
$ws.Range("G37").Value = 24
$ws.Range("H37").Value = 1788

$ws.Range("A38").Value = "Republica Dominicana"
$ws.Range("B38").Value = 74295
$ws.Range("C38").Value = 0
$ws.Range("D38").Value = 38824
$ws.Range("E38").Value = 34258
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = 0
$ws.Range("H38").Value = 1213

# 3) Armenia (row 54) updated figures
$ws.Range("B54").Value = 39586
$ws.Range("C54").Value = 288
$ws.Range("D54").Value = 30850
$ws.Range("E54").Value = 7968

# 4) Letonia (row 141) updated figures
$ws.Range("B141").Value = 1257
$ws.Range("C141").Value = 8
$ws.Range("E141").Value = 155

# 5) Georgia (row 145) updated figures
$ws.Range("B145").Value = 1197
$ws.Range("C145").Value = 15
$ws.Range("D145").Value = 974
$ws.Range("E145").Value = 206

# 6) Taiwan (row 166) updated figures
$ws.Range("D166").Value = 443
$ws.Range("E166").Value = 26
